# Fix Training Data Issue (#48)
# The "Date" column (BF) held a bogus value "4-7-2007-08" (a mangled
# concat of month-day and season) for every team row. It should be the
# actual game date in ISO form: 2008-04-07.
#
# Rows 2-31 (one per NBA team) all carry the same bad value in column BF
# (column 58). We rewrite each of them to the corrected date string.
#
# NumberFormat is forced to Text ("@") before the write so Excel's
# automatic type inference does not reinterpret "2008-04-07" as a date
# serial number, then ClearFormats() restores the cell to its original
# (default/no explicit style) formatting so only the textual content
# changes, matching the source data fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = 58  # column BF
$firstRow = 2
$lastRow = 31
$correctedDate = "2008-04-07"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    $cell.NumberFormat = "@"
    $cell.Value = $correctedDate
    $cell.ClearFormats()
}
